# "One more result for test": a new student ("Amir Zornić") is inserted
# as a new row 19 in the results table. The existing row 19
# ("Hasan Mujanović") shifts down to row 20, and the "UKUPNO" total
# formula in column G is (re)applied to both the new row and the row that
# got pushed down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the current last row of data (row 19, "Hasan Mujanović") down to
# row 20 by copying its values first...
$ws.Range("A20").Value = $ws.Range("A19").Value()
$ws.Range("B20").Value = $ws.Range("B19").Value()
$ws.Range("C20").Value = $ws.Range("C19").Value()
$ws.Range("D20").Value = $ws.Range("D19").Value()
$ws.Range("E20").Value = $ws.Range("E19").Value()
$ws.Range("F20").Value = $ws.Range("F19").Value()
$ws.Range("G20").Formula = "=C20+D20+E20+F20"

# ...then cloning row 19's look (borders/fills/number formats) onto row 20.
$ws.Range("A19:G19").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)

# Write the new student's data into row 19.
$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "Amir Zornić"
$ws.Range("C19").Value = 30
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = 30
$ws.Range("F19").Value = 20
$ws.Range("G19").Formula = "=C19+D19+E19+F19"
